$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 14-16: filled in normally, column by column across each row
$rows1 = @(
    @(14, "013", 804, 571, 1040, 703, "after match - you win"),
    @(15, "014", 994, 594, 1167, 682, "after match - break down"),
    @(16, "015", 656, 488, 1191, 589, "after match - rank up")
)
foreach ($row in $rows1) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

# Rows 17-19: column A filled first for all three rows (016, 017, 018)
$ws.Cells.Item(17, 1).Value = "016"
$ws.Cells.Item(18, 1).Value = "017"
$ws.Cells.Item(19, 1).Value = "018"

# Then remaining numeric columns for rows 17-19
$numericRows = @(
    @(17, 716, 147, 1124, 223),
    @(18, 639, 826, 823,  901),
    @(19, 1012, 827, 1204, 902)
)
foreach ($row in $numericRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Then column F filled bottom-to-top: row 19, 18, 17
$ws.Cells.Item(19, 6).Value = "recover energy - restore button"
$ws.Cells.Item(18, 6).Value = "recover energy - cancel button"
$ws.Cells.Item(17, 6).Value = "recover energy title"

$ws.Range("F17").Select()
